# Applies the Betfair odds updates for 2025-12-17 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Serbian Super League: FK Radnicki 1923 vs Cukaricki
$ws.Range("F2").Value = 2.24
$ws.Range("H2").Value = 2.58
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 3.9
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.18
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 2.42
$ws.Range("T2").Value = 1.57
$ws.Range("U2").Value = 2.18

# Row 3 - Swiss Super League: Young Boys vs Grasshoppers Zurich
$ws.Range("F3").Value = 1.59
$ws.Range("G3").Value = 1.62
$ws.Range("H3").Value = 5.3
$ws.Range("O3").Value = 1.16
$ws.Range("Q3").Value = 1.48
$ws.Range("T3").Value = 1.51
$ws.Range("U3").Value = 2.46
$ws.Range("V3").Value = 1.2
$ws.Range("W3").Value = 2.58
$ws.Range("X3").Value = 30

# Row 4 - Swiss Super League: FC Zurich vs Lugano
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2.44
$ws.Range("L4").Value = 1.33
$ws.Range("Q4").Value = 1.68
$ws.Range("U4").Value = 2.16
$ws.Range("AI4").Value = 38
$ws.Range("AL4").Value = 42
$ws.Range("AN4").Value = 25

# Row 5 - Swiss Super League: Luzern vs FC Basel
$ws.Range("I5").Value = 2.14
$ws.Range("L5").Value = 1.27
$ws.Range("R5").Value = 1.64
$ws.Range("V5").Value = 1.88
$ws.Range("W5").Value = 1.37

# Row 6 - Scottish Premiership: Dundee Utd vs Celtic
$ws.Range("I6").Value = 1.53
$ws.Range("T6").Value = 1.8
$ws.Range("V6").Value = 2.9
